$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Merge old row 93 ("General" header, column B only) into old row 94
#    (G1.FB / text), producing a single combined row, then delete the
#    now-redundant row 93 and re-insert 3 fresh rows above the merged
#    row so it lands on row 96 (matching the target layout) while
#    leaving rows 93/94 free for two brand new backlog items and row 95
#    as an empty separator row (as in the rest of the sheet).
# ------------------------------------------------------------------
$ws.Range("B94").Value = $ws.Range("B93").Value2
$ws.Rows(93).Delete()
$ws.Rows("93:95").Insert()

# The Insert() above clones formatting from the row above (row 92), so
# strip the stray formatted-but-empty cells it leaves behind.
$ws.Range("D93").Clear()
$ws.Rows(95).Clear()

# ------------------------------------------------------------------
# 2) Row-height tweaks: rows 92-94 become the shorter "17.25" custom
#    height used by the newer backlog rows.
# ------------------------------------------------------------------
$ws.Rows(92).RowHeight = 17.25
$ws.Rows(93).RowHeight = 17.25
$ws.Rows(94).RowHeight = 17.25

# Reuse the existing date-format style (style index already used by
# E92) for the two new date cells instead of synthesizing a new
# numFmt, by copy/pasting formats only.
$ws.Range("E92").Copy()
$ws.Range("E93:E94").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) New row 93: S23 - Modify the Tab item style ...
# ------------------------------------------------------------------
$ws.Range("A93").Value = "S23"
$ws.Range("C93").Value = "Modify the Tab item style to be fat finger compatible and more colorful on selection."

# ------------------------------------------------------------------
# 4) New row 94: S24 - Remember and reassert window placement ...
# ------------------------------------------------------------------
$ws.Range("A94").Value = "S24"
$ws.Range("C94").Value = "Remember and reassert window placement between sessions"
$ws.Range("D94").Value = "Done"
$ws.Range("E94").Value = 43555

# ------------------------------------------------------------------
# 5) Append the new "Toolbar" section at the bottom of the sheet
#    (rows 140-144, leaving row 139 blank as a separator like
#    elsewhere in the sheet).
# ------------------------------------------------------------------
$ws.Range("A140").Value = "T1"
$ws.Range("B140").Value = "Toolbar"

$ws.Range("A141").Value = "T2"
$ws.Range("C141").Value = "Add buttons that show and hide depending upon the currently active tab control"
$ws.Range("D141").Value = "Done"
$ws.Range("E141").Value = 43554
$ws.Range("E92").Copy()
$ws.Range("E141").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E141").Value = 43554

$ws.Range("A142").Value = "T3"
$ws.Range("C142").Value = "Ask G for an image for Copy Inventory button"

$ws.Range("A143").Value = "T4"
$ws.Range("C143").Value = "Ask G for an image for Copy Inspection button"

$ws.Range("A144").Value = "T5"
$ws.Range("C144").Value = "Ask G for an image for Copy Sections button"

# ------------------------------------------------------------------
# 6) View-state cosmetics captured by the diff: scroll position,
#    selection and the workbook window placement.
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 77
$ws.Range("D93").Select()

$excel.ActiveWindow.WindowState = -4143
$excel.Left = 2100
$excel.Top = 1155
